# TC01_INS_Filter_Doc-CCG.xlsx edit:
#  - The "startup" sheet's StatQuery table is updated so that the
#    former "ProjectsTab" row becomes a "GrantsTab" row, and all five
#    rows' StatQuery (column C) text is replaced with the new combined
#    Neo4j query (which now also returns Grants / uses Project IDs etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# New StatQuery text shared by every row (column C) of the table.
$newQuery = "MATCH (p:program)<--(pr:project)`n" +
            "where pr.lead_doc='CCG'`n" +
            "OPTIONAL MATCH (pr)<--(pub:publication)`n" +
            "OPTIONAL MATCH (ct:clinical_trial)`n" +
            "WHERE EXISTS((pr)<--(pub)<--(ct)) OR EXISTS((pr)<--(ct))`n" +
            "OPTIONAL MATCH (pr)<--(pat)`n" +
            "WHERE pat:patent_application OR pat:granted_patent`n" +
            "OPTIONAL MATCH (pr)<-[*1..2]-(dt)`n" +
            "WHERE dt:sra OR dt:dbgap OR dt:geo`n" +
            "WITH p, pr, pub, ct, pat, dt`n" +
            "RETURN`n" +
            "COUNT(DISTINCT p.program_id) AS Programs,`n" +
            "COUNT(DISTINCT pr.queried_project_id) AS Projects,`n" +
            "COUNT(DISTINCT pr.project_id) AS Grants,`n" +
            "COUNT(DISTINCT pub.publication_id) AS Publications,`n" +
            "COUNT(DISTINCT dt.accession) AS Datasets,`n" +
            "COUNT(DISTINCT ct.clinical_trial_id) AS ``Clinical Trials```,`n" +
            "COUNT(DISTINCT pat.patent_id) AS Patents"

# Update column C (StatQuery) for every data row first, so the new
# shared string is created before the "GrantsTab" label below.
$ws.Range("C2").Value2 = $newQuery
$ws.Range("C3").Value2 = $newQuery
$ws.Range("C4").Value2 = $newQuery
$ws.Range("C5").Value2 = $newQuery
$ws.Range("C6").Value2 = $newQuery

# Row 2 used to describe "ProjectsTab"; it now describes "GrantsTab"
# (its query, in column B, already targeted grants/projects and is
# left untouched).
$ws.Range("A2").Value2 = "GrantsTab"

# The longer query text now needs a taller row for row 2.
$ws.Rows.Item(2).RowHeight = 270

# Reflect the final selection/scroll position left in the sheet.
$ws.Activate()
$ws.Range("A6").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
